$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44162
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 100
$ws.Cells.Item(2, 14).Value = 4000
$ws.Cells.Item(2, 15).Value = 4000
$ws.Cells.Item(2, 16).Value = 4000
$ws.Cells.Item(2, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 2000
$ws.Cells.Item(2, 20).Value = 2

# Row 3
$ws.Cells.Item(3, 4).Value = 44264
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 110
$ws.Cells.Item(3, 14).Value = 3500
$ws.Cells.Item(3, 15).Value = 4000
$ws.Cells.Item(3, 16).Value = 3727
$ws.Cells.Item(3, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(3, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(3, 19).Value = 1864
$ws.Cells.Item(3, 20).Value = 2

# Row 4
$ws.Cells.Item(4, 4).Value = 44516
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 30
$ws.Cells.Item(4, 14).Value = 5000
$ws.Cells.Item(4, 15).Value = 5000
$ws.Cells.Item(4, 16).Value = 5000
$ws.Cells.Item(4, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(4, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(4, 19).Value = 2500
$ws.Cells.Item(4, 20).Value = 2

# Row 5
$ws.Cells.Item(5, 4).Value = 44187
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 110
$ws.Cells.Item(5, 14).Value = 2600
$ws.Cells.Item(5, 15).Value = 3000
$ws.Cells.Item(5, 16).Value = 2782
$ws.Cells.Item(5, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(5, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(5, 19).Value = 1391
$ws.Cells.Item(5, 20).Value = 2

# Row 6
$ws.Cells.Item(6, 4).Value = 44526
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 250
$ws.Cells.Item(6, 14).Value = 4000
$ws.Cells.Item(6, 15).Value = 4000
$ws.Cells.Item(6, 16).Value = 4000
$ws.Cells.Item(6, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(6, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(6, 19).Value = 2000
$ws.Cells.Item(6, 20).Value = 2

# Row 7
$ws.Cells.Item(7, 4).Value = 44235
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 60
$ws.Cells.Item(7, 14).Value = 3000
$ws.Cells.Item(7, 15).Value = 3000
$ws.Cells.Item(7, 16).Value = 3000
$ws.Cells.Item(7, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(7, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(7, 19).Value = 1500
$ws.Cells.Item(7, 20).Value = 2

# Row 8
$ws.Cells.Item(8, 4).Value = 44169
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 400
$ws.Cells.Item(8, 14).Value = 3600
$ws.Cells.Item(8, 15).Value = 3600
$ws.Cells.Item(8, 16).Value = 3600
$ws.Cells.Item(8, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(8, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(8, 19).Value = 1800
$ws.Cells.Item(8, 20).Value = 2

# Row 9
$ws.Cells.Item(9, 4).Value = 44524
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 180
$ws.Cells.Item(9, 14).Value = 4000
$ws.Cells.Item(9, 15).Value = 4000
$ws.Cells.Item(9, 16).Value = 4000
$ws.Cells.Item(9, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(9, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(9, 19).Value = 2000
$ws.Cells.Item(9, 20).Value = 2

# Row 10
$ws.Cells.Item(10, 4).Value = 44202
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 30
$ws.Cells.Item(10, 14).Value = 3000
$ws.Cells.Item(10, 15).Value = 3000
$ws.Cells.Item(10, 16).Value = 3000
$ws.Cells.Item(10, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(10, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(10, 19).Value = 1500
$ws.Cells.Item(10, 20).Value = 2

# Row 11
$ws.Cells.Item(11, 4).Value = 44202
$ws.Cells.Item(11, 12).Value = 'Segunda'
$ws.Cells.Item(11, 13).Value = 20
$ws.Cells.Item(11, 14).Value = 2600
$ws.Cells.Item(11, 15).Value = 2600
$ws.Cells.Item(11, 16).Value = 2600
$ws.Cells.Item(11, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(11, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(11, 19).Value = 1300
$ws.Cells.Item(11, 20).Value = 2

# Row 12
$ws.Cells.Item(12, 4).Value = 44511
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 50
$ws.Cells.Item(12, 14).Value = 6400
$ws.Cells.Item(12, 15).Value = 6400
$ws.Cells.Item(12, 16).Value = 6400
$ws.Cells.Item(12, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(12, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(12, 19).Value = 3200
$ws.Cells.Item(12, 20).Value = 2

# Row 13
$ws.Cells.Item(13, 4).Value = 44517
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 20
$ws.Cells.Item(13, 14).Value = 5000
$ws.Cells.Item(13, 15).Value = 5000
$ws.Cells.Item(13, 16).Value = 5000
$ws.Cells.Item(13, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(13, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(13, 19).Value = 2500
$ws.Cells.Item(13, 20).Value = 2

# Row 14
$ws.Cells.Item(14, 4).Value = 44265
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 70
$ws.Cells.Item(14, 14).Value = 3600
$ws.Cells.Item(14, 15).Value = 3800
$ws.Cells.Item(14, 16).Value = 3714
$ws.Cells.Item(14, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(14, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(14, 19).Value = 1857
$ws.Cells.Item(14, 20).Value = 2

# Row 15
$ws.Cells.Item(15, 4).Value = 44172
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 300
$ws.Cells.Item(15, 14).Value = 3400
$ws.Cells.Item(15, 15).Value = 3600
$ws.Cells.Item(15, 16).Value = 3467
$ws.Cells.Item(15, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(15, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(15, 19).Value = 1734
$ws.Cells.Item(15, 20).Value = 2

# Row 16
$ws.Cells.Item(16, 4).Value = 44186
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 3000
$ws.Cells.Item(16, 15).Value = 3000
$ws.Cells.Item(16, 16).Value = 3000
$ws.Cells.Item(16, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(16, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(16, 19).Value = 1500
$ws.Cells.Item(16, 20).Value = 2

# Row 17
$ws.Cells.Item(17, 4).Value = 44232
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 60
$ws.Cells.Item(17, 14).Value = 3000
$ws.Cells.Item(17, 15).Value = 3000
$ws.Cells.Item(17, 16).Value = 3000
$ws.Cells.Item(17, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(17, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(17, 19).Value = 1500
$ws.Cells.Item(17, 20).Value = 2

# Row 18
$ws.Cells.Item(18, 4).Value = 44176
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 150
$ws.Cells.Item(18, 14).Value = 3500
$ws.Cells.Item(18, 15).Value = 3500
$ws.Cells.Item(18, 16).Value = 3500
$ws.Cells.Item(18, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(18, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(18, 19).Value = 2333
$ws.Cells.Item(18, 20).Value = 1.5

# Row 19
$ws.Cells.Item(19, 4).Value = 44165
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 400
$ws.Cells.Item(19, 14).Value = 3400
$ws.Cells.Item(19, 15).Value = 3400
$ws.Cells.Item(19, 16).Value = 3400
$ws.Cells.Item(19, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(19, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 19).Value = 1700
$ws.Cells.Item(19, 20).Value = 2

# Row 20
$ws.Cells.Item(20, 4).Value = 44530
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 100
$ws.Cells.Item(20, 14).Value = 4000
$ws.Cells.Item(20, 15).Value = 4000
$ws.Cells.Item(20, 16).Value = 4000
$ws.Cells.Item(20, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(20, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(20, 19).Value = 2000
$ws.Cells.Item(20, 20).Value = 2

# Row 21
$ws.Cells.Item(21, 4).Value = 44530
$ws.Cells.Item(21, 12).Value = 'Segunda'
$ws.Cells.Item(21, 13).Value = 150
$ws.Cells.Item(21, 14).Value = 3600
$ws.Cells.Item(21, 15).Value = 3600
$ws.Cells.Item(21, 16).Value = 3600
$ws.Cells.Item(21, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(21, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(21, 19).Value = 1800
$ws.Cells.Item(21, 20).Value = 2

# Row 22
$ws.Cells.Item(22, 4).Value = 44529
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 100
$ws.Cells.Item(22, 14).Value = 3800
$ws.Cells.Item(22, 15).Value = 3800
$ws.Cells.Item(22, 16).Value = 3800
$ws.Cells.Item(22, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(22, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(22, 19).Value = 1900
$ws.Cells.Item(22, 20).Value = 2

# Row 23
$ws.Cells.Item(23, 4).Value = 44210
$ws.Cells.Item(23, 12).Value = 'Segunda'
$ws.Cells.Item(23, 13).Value = 150
$ws.Cells.Item(23, 14).Value = 2700
$ws.Cells.Item(23, 15).Value = 2700
$ws.Cells.Item(23, 16).Value = 2700
$ws.Cells.Item(23, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(23, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(23, 19).Value = 1350
$ws.Cells.Item(23, 20).Value = 2

# Row 24
$ws.Cells.Item(24, 4).Value = 44518
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 20
$ws.Cells.Item(24, 14).Value = 5000
$ws.Cells.Item(24, 15).Value = 5000
$ws.Cells.Item(24, 16).Value = 5000
$ws.Cells.Item(24, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(24, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(24, 19).Value = 2500
$ws.Cells.Item(24, 20).Value = 2

# Row 25
$ws.Cells.Item(25, 4).Value = 44519
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 180
$ws.Cells.Item(25, 14).Value = 4000
$ws.Cells.Item(25, 15).Value = 4000
$ws.Cells.Item(25, 16).Value = 4000
$ws.Cells.Item(25, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(25, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(25, 19).Value = 2000
$ws.Cells.Item(25, 20).Value = 2

# Row 26
$ws.Cells.Item(26, 4).Value = 44204
$ws.Cells.Item(26, 12).Value = 'Primera'
$ws.Cells.Item(26, 13).Value = 50
$ws.Cells.Item(26, 14).Value = 3000
$ws.Cells.Item(26, 15).Value = 3000
$ws.Cells.Item(26, 16).Value = 3000
$ws.Cells.Item(26, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(26, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(26, 19).Value = 1500
$ws.Cells.Item(26, 20).Value = 2

# Row 27
$ws.Cells.Item(27, 4).Value = 44204
$ws.Cells.Item(27, 12).Value = 'Segunda'
$ws.Cells.Item(27, 13).Value = 140
$ws.Cells.Item(27, 14).Value = 2400
$ws.Cells.Item(27, 15).Value = 2400
$ws.Cells.Item(27, 16).Value = 2400
$ws.Cells.Item(27, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(27, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(27, 19).Value = 1200
$ws.Cells.Item(27, 20).Value = 2

# Row 28
$ws.Cells.Item(28, 4).Value = 44200
$ws.Cells.Item(28, 12).Value = 'Segunda'
$ws.Cells.Item(28, 13).Value = 50
$ws.Cells.Item(28, 14).Value = 2600
$ws.Cells.Item(28, 15).Value = 2600
$ws.Cells.Item(28, 16).Value = 2600
$ws.Cells.Item(28, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(28, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(28, 19).Value = 1300
$ws.Cells.Item(28, 20).Value = 2

# Row 29
$ws.Cells.Item(29, 4).Value = 44166
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 1500
$ws.Cells.Item(29, 14).Value = 3600
$ws.Cells.Item(29, 15).Value = 3600
$ws.Cells.Item(29, 16).Value = 3600
$ws.Cells.Item(29, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(29, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(29, 19).Value = 1800
$ws.Cells.Item(29, 20).Value = 2

# Row 30
$ws.Cells.Item(30, 4).Value = 44522
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 30
$ws.Cells.Item(30, 14).Value = 5000
$ws.Cells.Item(30, 15).Value = 5000
$ws.Cells.Item(30, 16).Value = 5000
$ws.Cells.Item(30, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(30, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(30, 19).Value = 2500
$ws.Cells.Item(30, 20).Value = 2

# Row 31
$ws.Cells.Item(31, 4).Value = 44525
$ws.Cells.Item(31, 12).Value = 'Primera'
$ws.Cells.Item(31, 13).Value = 150
$ws.Cells.Item(31, 14).Value = 4000
$ws.Cells.Item(31, 15).Value = 4000
$ws.Cells.Item(31, 16).Value = 4000
$ws.Cells.Item(31, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(31, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(31, 19).Value = 2000
$ws.Cells.Item(31, 20).Value = 2

# Row 32
$ws.Cells.Item(32, 4).Value = 44512
$ws.Cells.Item(32, 12).Value = 'Primera'
$ws.Cells.Item(32, 13).Value = 30
$ws.Cells.Item(32, 14).Value = 6000
$ws.Cells.Item(32, 15).Value = 6000
$ws.Cells.Item(32, 16).Value = 6000
$ws.Cells.Item(32, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(32, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(32, 19).Value = 3000
$ws.Cells.Item(32, 20).Value = 2

# Row 33
$ws.Cells.Item(33, 4).Value = 44211
$ws.Cells.Item(33, 12).Value = 'Primera'
$ws.Cells.Item(33, 13).Value = 40
$ws.Cells.Item(33, 14).Value = 2800
$ws.Cells.Item(33, 15).Value = 2800
$ws.Cells.Item(33, 16).Value = 2800
$ws.Cells.Item(33, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(33, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(33, 19).Value = 1400
$ws.Cells.Item(33, 20).Value = 2

# Row 34
$ws.Cells.Item(34, 4).Value = 44211
$ws.Cells.Item(34, 12).Value = 'Segunda'
$ws.Cells.Item(34, 13).Value = 30
$ws.Cells.Item(34, 14).Value = 2600
$ws.Cells.Item(34, 15).Value = 2600
$ws.Cells.Item(34, 16).Value = 2600
$ws.Cells.Item(34, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(34, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(34, 19).Value = 1300
$ws.Cells.Item(34, 20).Value = 2

# Row 35
$ws.Cells.Item(35, 4).Value = 44167
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 500
$ws.Cells.Item(35, 14).Value = 3600
$ws.Cells.Item(35, 15).Value = 3600
$ws.Cells.Item(35, 16).Value = 3600
$ws.Cells.Item(35, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(35, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(35, 19).Value = 1800
$ws.Cells.Item(35, 20).Value = 2
